$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").Value = "Коньяк"
$ws.Range("B34").Value = 0

$ws.Range("B35").Select()
